# "add brands & change db"
#
# The orderList sheet now records the shipping option chosen for each
# order line (a new "shippingId" column between productId and quantity).
# Since that info moved down to the per-line table, the single shippingId
# column on the order sheet (one per order) is no longer needed and is
# removed.

$wb = $excel.ActiveWorkbook

# --- orderList: insert a "shippingId" column (C) with per-row values ---
$wsOrderList = $wb.Worksheets.Item("orderList")
$wsOrderList.Columns("C:C").Insert()
$wsOrderList.Range("C1").Value = "shippingId"
$wsOrderList.Range("C2").Value = 1
$wsOrderList.Range("C3").Value = 1
$wsOrderList.Range("C4").Value = 2
$wsOrderList.Range("C5").Value = 1
$wsOrderList.Range("C6").Value = 2
$wsOrderList.Range("C7").Value = 3

# --- order: drop the old shippingId column (C) ---
$wsOrder = $wb.Worksheets.Item("order")
$wsOrder.Columns("C:C").Delete()

# --- restore / update the selections + active sheet ---
$wsOrder.Range("C1:I4").Select()
$wsOrderList.Range("E11").Select()
